$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 71,3
$arr[0,0] = 22.3645970149253
$arr[0,1] = 14.808
$arr[0,2] = 30.19199999999999
$arr[1,0] = 2.530228855721378
$arr[1,1] = 0.9559999999999998
$arr[1,2] = 4.359999999999998
$arr[2,0] = 4.446407960198998
$arr[2,1] = 1.956
$arr[2,2] = 7.280000000000001
$arr[3,0] = 3.380179104477605
$arr[3,1] = 1.388
$arr[3,2] = 5.888000000000004
$arr[4,0] = 11.09305472636814
$arr[4,1] = 6.435999999999999
$arr[4,2] = 16.33599999999999
$arr[5,0] = 17.76332338308457
$arr[5,1] = 11.36
$arr[5,2] = 25.728
$arr[6,0] = 27.91032835820883
$arr[6,1] = 20.36800000000001
$arr[6,2] = 35.37199999999999
$arr[7,0] = 9.832099502487553
$arr[7,1] = 5.399999999999998
$arr[7,2] = 15.056
$arr[8,0] = 30.03215920397993
$arr[8,1] = 22.38399999999999
$arr[8,2] = 37.188
$arr[9,0] = 27.09962189054717
$arr[9,1] = 18.764
$arr[9,2] = 35.40399999999998
$arr[10,0] = 19.6647960199005
$arr[10,1] = 11.616
$arr[10,2] = 28.45999999999999
$arr[11,0] = 35.20533333333324
$arr[11,1] = 25.58399999999999
$arr[11,2] = 44.60800000000003
$arr[12,0] = 30.95683582089534
$arr[12,1] = 22.02800000000001
$arr[12,2] = 40.37200000000001
$arr[13,0] = 4.165532338308449
$arr[13,1] = 1.835999999999999
$arr[13,2] = 7.083999999999998
$arr[14,0] = 11.00141293532337
$arr[14,1] = 6.328
$arr[14,2] = 16.224
$arr[15,0] = 19.44220895522386
$arr[15,1] = 12.364
$arr[15,2] = 28.10399999999999
$arr[16,0] = 4.432278606965167
$arr[16,1] = 1.944
$arr[16,2] = 7.248000000000001
$arr[17,0] = 30.89245771144274
$arr[17,1] = 21.096
$arr[17,2] = 40.90400000000002
$arr[18,0] = 13.22837810945273
$arr[18,1] = 7.756000000000005
$arr[18,2] = 20.128
$arr[19,0] = 7.319223880597013
$arr[19,1] = 3.639999999999997
$arr[19,2] = 11.812
$arr[20,0] = 32.73251741293522
$arr[20,1] = 24.04
$arr[20,2] = 41.584
$arr[21,0] = 4.681691542288546
$arr[21,1] = 1.952
$arr[21,2] = 7.768000000000004
$arr[22,0] = 11.08803980099502
$arr[22,1] = 6.340000000000001
$arr[22,2] = 16.416
$arr[23,0] = 8.188656716417899
$arr[23,1] = 4.684000000000003
$arr[23,2] = 12.784
$arr[24,0] = 7.302129353233826
$arr[24,1] = 3.891999999999999
$arr[24,2] = 11.368
$arr[25,0] = 38.4543681592039
$arr[25,1] = 26.92000000000001
$arr[25,2] = 48.84000000000001
$arr[26,0] = 20.09554228855719
$arr[26,1] = 12.96
$arr[26,2] = 26.996
$arr[27,0] = 4.971422885572127
$arr[27,1] = 2.104000000000001
$arr[27,2] = 8.176000000000005
$arr[28,0] = 17.80356218905472
$arr[28,1] = 11.092
$arr[28,2] = 26.51199999999999
$arr[29,0] = 29.02171144278596
$arr[29,1] = 21.01999999999999
$arr[29,2] = 37.39600000000001
$arr[30,0] = 28.21785074626855
$arr[30,1] = 20.09600000000001
$arr[30,2] = 36.30000000000003
$arr[31,0] = 31.10461691542272
$arr[31,1] = 22.21600000000002
$arr[31,2] = 41.05199999999997
$arr[32,0] = 26.91713432835814
$arr[32,1] = 17.02399999999999
$arr[32,2] = 36.48400000000001
$arr[33,0] = 12.01729353233831
$arr[33,1] = 6.812000000000003
$arr[33,2] = 18.44400000000001
$arr[34,0] = 17.73002985074625
$arr[34,1] = 11.42399999999999
$arr[34,2] = 25.968
$arr[35,0] = 8.230089552238793
$arr[35,1] = 4.256000000000002
$arr[35,2] = 13.068
$arr[36,0] = 34.37277611940294
$arr[36,1] = 25.32
$arr[36,2] = 43.84400000000002
$arr[37,0] = 21.30614925373134
$arr[37,1] = 13.736
$arr[37,2] = 30.68
$arr[38,0] = 5.781990049751236
$arr[38,1] = 2.852
$arr[38,2] = 9.823999999999995
$arr[39,0] = 15.63026865671641
$arr[39,1] = 10.32799999999999
$arr[39,2] = 22.068
$arr[40,0] = 29.659303482587
$arr[40,1] = 19.9
$arr[40,2] = 39.70800000000001
$arr[41,0] = 17.75420895522387
$arr[41,1] = 11.064
$arr[41,2] = 26.42799999999999
$arr[42,0] = 19.76847761194028
$arr[42,1] = 12.708
$arr[42,2] = 26.62
$arr[43,0] = 32.56483582089545
$arr[43,1] = 23.07599999999998
$arr[43,2] = 41.76799999999997
$arr[44,0] = 26.17456716417904
$arr[44,1] = 17.968
$arr[44,2] = 34.23599999999998
$arr[45,0] = 34.68195024875612
$arr[45,1] = 24.96399999999998
$arr[45,2] = 43.56799999999998
$arr[46,0] = 19.66658706467662
$arr[46,1] = 11.592
$arr[46,2] = 28.29199999999998
$arr[47,0] = 4.957014925373131
$arr[47,1] = 2.547999999999998
$arr[47,2] = 8.036
$arr[48,0] = 8.194547263681583
$arr[48,1] = 4.100000000000001
$arr[48,2] = 13.04399999999999
$arr[49,0] = 2.518985074626852
$arr[49,1] = 0.9479999999999998
$arr[49,2] = 4.327999999999997
$arr[50,0] = 23.25751243781087
$arr[50,1] = 16.29999999999999
$arr[50,2] = 30.17200000000001
$arr[51,0] = 20.9857910447761
$arr[51,1] = 13.492
$arr[51,2] = 29.824
$arr[52,0] = 4.166507462686559
$arr[52,1] = 1.839999999999999
$arr[52,2] = 7.083999999999998
$arr[53,0] = 7.30947263681592
$arr[53,1] = 3.655999999999997
$arr[53,2] = 11.8
$arr[54,0] = 28.57586069651737
$arr[54,1] = 19.73200000000001
$arr[54,2] = 36.74799999999996
$arr[55,0] = 28.21297512437799
$arr[55,1] = 20.11200000000001
$arr[55,2] = 36.26400000000002
$arr[56,0] = 8.141810945273622
$arr[56,1] = 4.256000000000003
$arr[56,2] = 13.212
$arr[57,0] = 25.29745273631822
$arr[57,1] = 16.72000000000001
$arr[57,2] = 34.08399999999999
$arr[58,0] = 4.718646766169143
$arr[58,1] = 2.008
$arr[58,2] = 7.683999999999998
$arr[59,0] = 33.07988059701481
$arr[59,1] = 22.86399999999998
$arr[59,2] = 43.15199999999998
$arr[60,0] = 27.09962189054717
$arr[60,1] = 18.764
$arr[60,2] = 35.40399999999998
$arr[61,0] = 36.61056716417895
$arr[61,1] = 26.01599999999999
$arr[61,2] = 46.188
$arr[62,0] = 4.824537313432824
$arr[62,1] = 2.048
$arr[62,2] = 7.968000000000003
$arr[63,0] = 25.83908457711442
$arr[63,1] = 18.624
$arr[63,2] = 33.21600000000002
$arr[64,0] = 25.8221293532337
$arr[64,1] = 16.69599999999999
$arr[64,2] = 35.71600000000002
$arr[65,0] = 23.97028855721383
$arr[65,1] = 16.104
$arr[65,2] = 34.49600000000002
$arr[66,0] = 6.586427860696502
$arr[66,1] = 3.531999999999999
$arr[66,2] = 10.704
$arr[67,0] = 34.9284776119402
$arr[67,1] = 25.17999999999998
$arr[67,2] = 44.97600000000002
$arr[68,0] = 5.844776119402974
$arr[68,1] = 2.924
$arr[68,2] = 9.555999999999994
$arr[69,0] = 11.04045771144278
$arr[69,1] = 6.383999999999999
$arr[69,2] = 16.964
$arr[70,0] = 14.88131343283582
$arr[70,1] = 9.159999999999995
$arr[70,2] = 21.43999999999999
$ws.Range("A2:C72").Value = $arr
